# Update the names on the "employees-valid" sheet (sheet 1) to the new
# Lord-of-the-Rings-themed data set, then leave that sheet selected with
# C2 as the active cell (matching the author's final view state).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 2: Jack Sparrow -> Gandalf Thegrey (dependents unchanged: 10)
$ws1.Range("B2").Value = "Gandalf"
$ws1.Range("A2").Value = "Thegrey"

# Row 3: Hector Barbossa -> Frodo Baggins (dependents: 1 -> 32)
$ws1.Range("B3").Value = "Frodo"
$ws1.Range("A3").Value = "Baggins"
$ws1.Range("C3").Value = "32"

# Row 4: Davy Jones -> Sam Gamgee (dependents: 32 -> 1)
$ws1.Range("B4").Value = "Sam"
$ws1.Range("A4").Value = "Gamgee"
$ws1.Range("C4").Value = "1"

# Row 5: Joshamee Gibbs -> Gollum Smeagol (dependents unchanged: 0)
$ws1.Range("A5").Value = "Gollum"
$ws1.Range("B5").Value = "Smeagol"

# Make the "employees-valid" sheet the active/selected tab, with C2
# selected as the active cell (the "employees-updated" sheet, previously
# active, loses its tabSelected flag as a result).
$ws1.Activate()
$ws1.Range("C2").Select()
